$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3051.75
$ws.Range("I19").Value = 1537.6
$ws.Range("J19").Value = 4133.2856
$ws.Range("K19").Value = 1537.6
$ws.Range("L19").Value = 4133.2856
$ws.Range("M19").Value = -1362.6
$ws.Range("N19").Value = -4483.2856
$ws.Range("H28").Value = 713.84
$ws.Range("J28").Value = 1302.6
$ws.Range("L28").Value = 1302.6
$ws.Range("N28").Value = -2272.6
$ws.Range("H62").Value = 3508.0833
$ws.Range("I62").Value = 2928.4285
$ws.Range("K62").Value = 2928.4285
$ws.Range("M62").Value = -2304.4285
$ws.Range("H65").Value = 3508.0833
$ws.Range("I65").Value = 2928.4285
$ws.Range("K65").Value = 14642.1425
$ws.Range("M65").Value = -11522.1425
$ws.Range("H107").Value = 1835.6364
$ws.Range("I107").Value = 1600.5555
$ws.Range("K107").Value = 1600.5555
$ws.Range("M107").Value = 319.4445000000001
$ws.Range("H113").Value = 5260.5
$ws.Range("I113").Value = 5312.6
$ws.Range("K113").Value = 5312.6
$ws.Range("M113").Value = -2058.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3831.7046
$ws.Range("I61").Value = 707.1923
$ws.Range("J61").Value = 8344.888999999999
$ws.Range("K61").Value = 707.1923
$ws.Range("L61").Value = 8344.888999999999
$ws.Range("M61").Value = -495.1923
$ws.Range("N61").Value = -8768.888999999999
$ws.Range("H74").Value = 153602.34
$ws.Range("I74").Value = 194277.88
$ws.Range("J74").Value = 13497.777
$ws.Range("K74").Value = 194277.88
$ws.Range("L74").Value = 13497.777
$ws.Range("M74").Value = -193403.88
$ws.Range("N74").Value = -15245.777
$ws.Range("H77").Value = 153602.34
$ws.Range("I77").Value = 194277.88
$ws.Range("J77").Value = 13497.777
$ws.Range("K77").Value = 971389.4
$ws.Range("L77").Value = 67488.88499999999
$ws.Range("M77").Value = -967021.4
$ws.Range("N77").Value = -76224.88499999999
$ws.Range("H136").Value = 3831.7046
$ws.Range("I136").Value = 707.1923
$ws.Range("J136").Value = 8344.888999999999
$ws.Range("K136").Value = 2121.5769
$ws.Range("L136").Value = 25034.667
$ws.Range("M136").Value = 428.4231
$ws.Range("N136").Value = -30134.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 10015
$ws.Range("J21").Value = 10015
$ws.Range("L21").Value = 10015
$ws.Range("N21").Value = -10485
$ws.Range("H22").Value = 749.8889
$ws.Range("I22").Value = 425
$ws.Range("J22").Value = 842.7143
$ws.Range("K22").Value = 425
$ws.Range("L22").Value = 842.7143
$ws.Range("M22").Value = -75
$ws.Range("N22").Value = -1542.7143
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H59").Value = 49804
$ws.Range("J59").Value = 49804
$ws.Range("L59").Value = 49804
$ws.Range("N59").Value = -52094
$ws.Range("H94").Value = 1683.591
$ws.Range("I94").Value = 1458.875
$ws.Range("K94").Value = 1458.875
$ws.Range("M94").Value = -1007.875
$ws.Range("H107").Value = 2666.6667
$ws.Range("I107").Value = 4000
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 4000
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -2080
$ws.Range("N107").Value = -5840
$ws.Range("H132").Value = 31648.818
$ws.Range("I132").Value = 36718.68
$ws.Range("K132").Value = 110156.04
$ws.Range("M132").Value = -107626.04

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 4586.5
$ws.Range("I22").Value = 624.5
$ws.Range("K22").Value = 1873.5
$ws.Range("M22").Value = -1704.5
$ws.Range("H27").Value = 4586.5
$ws.Range("I27").Value = 624.5
$ws.Range("K27").Value = 1873.5
$ws.Range("M27").Value = -1771.5
$ws.Range("H41").Value = 2509.0908
$ws.Range("I41").Value = 450
$ws.Range("K41").Value = 1350
$ws.Range("M41").Value = -1012
$ws.Range("H44").Value = 4189.9
$ws.Range("J44").Value = 4322.1113
$ws.Range("L44").Value = 12966.3339
$ws.Range("N44").Value = -13762.3339
$ws.Range("H64").Value = 4862.0347
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 15000
$ws.Range("N64").Value = -15540
$ws.Range("H67").Value = 4862.0347
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 15000
$ws.Range("N67").Value = -16872
$ws.Range("H87").Value = 13313.083
$ws.Range("I87").Value = 9838
$ws.Range("J87").Value = 13809.523
$ws.Range("K87").Value = 29514
$ws.Range("L87").Value = 41428.569
$ws.Range("M87").Value = -28266
$ws.Range("N87").Value = -43924.569
$ws.Range("H90").Value = 13313.083
$ws.Range("I90").Value = 9838
$ws.Range("J90").Value = 13809.523
$ws.Range("K90").Value = 88542
$ws.Range("L90").Value = 124285.707
$ws.Range("M90").Value = -82302
$ws.Range("N90").Value = -136765.707
$ws.Range("H92").Value = 1413.7858
$ws.Range("I92").Value = 1823.125
$ws.Range("K92").Value = 5469.375
$ws.Range("M92").Value = -4221.375
$ws.Range("H98").Value = 512.94116
$ws.Range("I98").Value = 286.2
$ws.Range("J98").Value = 607.4167
$ws.Range("K98").Value = 858.5999999999999
$ws.Range("L98").Value = 1822.2501
$ws.Range("M98").Value = 639.4000000000001
$ws.Range("N98").Value = -4818.2501
$ws.Range("H99").Value = 11973.658
$ws.Range("J99").Value = 12541.639
$ws.Range("L99").Value = 37624.917
$ws.Range("N99").Value = -42116.917
$ws.Range("H117").Value = 2162.3333
$ws.Range("I117").Value = 2000
$ws.Range("K117").Value = 6000
$ws.Range("M117").Value = -2558
$ws.Range("H132").Value = 2542.2856
$ws.Range("I132").Value = 3274.5
$ws.Range("J132").Value = 1566
$ws.Range("K132").Value = 29470.5
$ws.Range("L132").Value = 14094
$ws.Range("M132").Value = -26940.5
$ws.Range("N132").Value = -19154
$ws.Range("H133").Value = 3707.3333
$ws.Range("I133").Value = 3117.125
$ws.Range("K133").Value = 9351.375
$ws.Range("M133").Value = -4291.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 7011.875
$ws.Range("I55").Value = 3519
$ws.Range("J55").Value = 12833.333
$ws.Range("K55").Value = 3519
$ws.Range("L55").Value = 12833.333
$ws.Range("M55").Value = -3192
$ws.Range("N55").Value = -13487.333
$ws.Range("H80").Value = 18269.4
$ws.Range("I80").Value = 2436
$ws.Range("J80").Value = 28825
$ws.Range("K80").Value = 2436
$ws.Range("L80").Value = 28825
$ws.Range("M80").Value = -1438
$ws.Range("N80").Value = -30821
$ws.Range("H83").Value = 18269.4
$ws.Range("I83").Value = 2436
$ws.Range("J83").Value = 28825
$ws.Range("K83").Value = 12180
$ws.Range("L83").Value = 144125
$ws.Range("M83").Value = -7188
$ws.Range("N83").Value = -154109
$ws.Range("H122").Value = 3043.8333
$ws.Range("I122").Value = 2782.25
$ws.Range("J122").Value = 3959.375
$ws.Range("K122").Value = 8346.75
$ws.Range("L122").Value = 11878.125
$ws.Range("M122").Value = -5896.75
$ws.Range("N122").Value = -16778.125
$ws.Range("H126").Value = 3570.3333
$ws.Range("I126").Value = 1893.8
$ws.Range("J126").Value = 5666
$ws.Range("K126").Value = 5681.4
$ws.Range("L126").Value = 16998
$ws.Range("M126").Value = -3211.4
$ws.Range("N126").Value = -21938
$ws.Range("H132").Value = 2793.8125
$ws.Range("I132").Value = 2399.0715
$ws.Range("J132").Value = 5557
$ws.Range("K132").Value = 7197.2145
$ws.Range("L132").Value = 16671
$ws.Range("M132").Value = -4667.2145
$ws.Range("N132").Value = -21731

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2566.2856
$ws.Range("I46").Value = 865.3684
$ws.Range("K46").Value = 865.3684
$ws.Range("M46").Value = -677.3684
$ws.Range("H82").Value = 1711.6364
$ws.Range("I82").Value = 1416.25
$ws.Range("J82").Value = 2499.3333
$ws.Range("K82").Value = 1416.25
$ws.Range("L82").Value = 2499.3333
$ws.Range("M82").Value = -1055.25
$ws.Range("N82").Value = -3221.3333
$ws.Range("H85").Value = 1711.6364
$ws.Range("I85").Value = 1416.25
$ws.Range("J85").Value = 2499.3333
$ws.Range("K85").Value = 1416.25
$ws.Range("L85").Value = 2499.3333
$ws.Range("M85").Value = -168.25
$ws.Range("N85").Value = -4995.3333
$ws.Range("H93").Value = 3193.182
$ws.Range("I93").Value = 2905
$ws.Range("J93").Value = 4490
$ws.Range("K93").Value = 2905
$ws.Range("L93").Value = 4490
$ws.Range("M93").Value = -1657
$ws.Range("N93").Value = -6986
$ws.Range("H122").Value = 6421.385
$ws.Range("I122").Value = 4866.6665
$ws.Range("J122").Value = 7754
$ws.Range("K122").Value = 14599.9995
$ws.Range("L122").Value = 23262
$ws.Range("M122").Value = -12149.9995
$ws.Range("N122").Value = -28162
$ws.Range("H136").Value = 3499.5945
$ws.Range("I136").Value = 3371.8076
$ws.Range("J136").Value = 3801.6365
$ws.Range("K136").Value = 10115.4228
$ws.Range("L136").Value = 11404.9095
$ws.Range("M136").Value = -7565.4228
$ws.Range("N136").Value = -16504.9095
